$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "time"
$ws.Range("D1").Value = "year"

$ws.Range("D1").Select()
